$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "Conclusi i file Astah degli scenari negativi." -> fix the typo
#    "Astah" -> "Asta". This also leaves the run split into three runs (as
#    happens in real Word when a word is edited in the middle and the
#    surrounding text keeps its own run identity): "Conclusi i file Asta" |
#    " " | "degli scenari negativi."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Astah", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Asta", 2) | Out-Null

$astaParagraph = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*Conclusi i file Asta*") {
        $astaParagraph = $d.Paragraphs.Item($i)
        break
    }
}

$base = $astaParagraph.Range.Start
$run1 = $d.Range($base, $base + 20)          # "Conclusi i file Asta"
$run2 = $d.Range($base + 20, $base + 21)     # " "
$run3 = $d.Range($base + 21, $base + 45)     # "degli scenari negativi."

# Nudge the font size away from and back to its original value on each
# sub-range so that Word materializes three separate runs instead of one.
foreach ($rr in @($run1, $run2, $run3)) {
    $rr.Font.Size = 20
    $rr.Font.Size = 14
}

# ---------------------------------------------------------------------------
# 2. Turn the trailing empty paragraph into a new list bullet that mirrors
#    the formatting of the preceding "Concluso il codice..." list item and
#    holds the new "scenari alternativi" sentence.
# ---------------------------------------------------------------------------
$codiceParagraph = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*Concluso il codice degli scenari negativi*") {
        $codiceParagraph = $d.Paragraphs.Item($i)
        break
    }
}

# Inserting a paragraph mark after an existing list item produces a new
# paragraph that inherits the same pStyle/numPr/rPr (list bullet, numId 3).
$codiceParagraph.Range.InsertParagraphAfter()
$newParagraph = $codiceParagraph.Next()

$newTextRange = $newParagraph.Range
$newTextRange.MoveEnd(1, -1) | Out-Null
$newTextRange.Text = "Conclusi i file Asta degli scenari alternativi (potrebbero variare se varia il codice)."

# Remove the now-redundant original trailing empty paragraph by merging its
# paragraph mark into the one we just wrote text into.
$trailingParagraph = $newParagraph.Next()
$mergeRange = $d.Range($newParagraph.Range.End - 1, $trailingParagraph.Range.End)
$mergeRange.Delete() | Out-Null
